$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Vostro / Dell" enrollment-device rows (ids 589 and 638);
# Excel shifts the remaining rows up and the now-unused shared strings
# are dropped automatically on save.
$ws.Rows.Item(6).Resize(2).Delete() | Out-Null

# Match the page setup used for printing (A4, portrait).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Leave the selection where the author left it after the edit.
$ws.Range("E16").Select() | Out-Null
